$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 261.33334
$ws.Range("I2").Value = 281.5
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 281.5
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = -168.5
$ws.Range("N2").Value = -326

$ws.Range("H70").Value = 2720.6667
$ws.Range("I70").Value = 2177.5
$ws.Range("J70").Value = 3399.625
$ws.Range("K70").Value = 6532.5
$ws.Range("L70").Value = 10198.875
$ws.Range("M70").Value = -6262.5
$ws.Range("N70").Value = -10738.875

$ws.Range("H73").Value = 2720.6667
$ws.Range("I73").Value = 2177.5
$ws.Range("J73").Value = 3399.625
$ws.Range("K73").Value = 6532.5
$ws.Range("L73").Value = 10198.875
$ws.Range("M73").Value = -5596.5
$ws.Range("N73").Value = -12070.875

$ws.Range("H86").Value = 5303.1816
$ws.Range("I86").Value = 4224.6924
$ws.Range("J86").Value = 6861
$ws.Range("K86").Value = 4224.6924
$ws.Range("L86").Value = 6861
$ws.Range("M86").Value = -3101.6924
$ws.Range("N86").Value = -9107

$ws.Range("H89").Value = 5303.1816
$ws.Range("I89").Value = 4224.6924
$ws.Range("J89").Value = 6861
$ws.Range("K89").Value = 21123.462
$ws.Range("L89").Value = 34305
$ws.Range("M89").Value = -15507.462
$ws.Range("N89").Value = -45537

$ws.Range("H92").Value = 1021.2727
$ws.Range("I92").Value = 594.26666
$ws.Range("K92").Value = 594.26666
$ws.Range("M92").Value = 653.73334

$ws.Range("H98").Value = 90909750
$ws.Range("I98").Value = 125000550
$ws.Range("K98").Value = 125000550
$ws.Range("M98").Value = -124999052

$ws.Range("H106").Value = 11138.167
$ws.Range("I106").Value = 1379.5
$ws.Range("K106").Value = 1379.5
$ws.Range("M106").Value = -748.5

$ws.Range("H115").Value = 6111
$ws.Range("I115").Value = 6538.875
$ws.Range("J115").Value = 4399.5
$ws.Range("K115").Value = 19616.625
$ws.Range("L115").Value = 13198.5
$ws.Range("M115").Value = -18049.625
$ws.Range("N115").Value = -16332.5

$ws.Range("H118").Value = 1137.8334
$ws.Range("I118").Value = 1165.4
$ws.Range("K118").Value = 3496.2
$ws.Range("M118").Value = -1839.2

$ws.Range("H121").Value = 1684998.4
$ws.Range("I121").Value = 2000
$ws.Range("J121").Value = 1783998.2
$ws.Range("K121").Value = 6000
$ws.Range("L121").Value = 5351994.6
$ws.Range("M121").Value = -4253
$ws.Range("N121").Value = -5355488.6

$ws.Range("H122").Value = 90909750
$ws.Range("I122").Value = 125000550
$ws.Range("K122").Value = 375001650
$ws.Range("M122").Value = -374999200

$ws.Range("H129").Value = 1707.1428
$ws.Range("J129").Value = 2010
$ws.Range("L129").Value = 6030
$ws.Range("N129").Value = -16030

$ws.Range("H138").Value = 3823.705
$ws.Range("J138").Value = 4034.0625
$ws.Range("L138").Value = 12102.1875
$ws.Range("N138").Value = -22382.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7474718
$ws.Range("I32").Value = 8340381
$ws.Range("K32").Value = 8340381
$ws.Range("M32").Value = -8340094

$ws.Range("H33").Value = 69950
$ws.Range("I33").Value = 69950
$ws.Range("K33").Value = 69950
$ws.Range("M33").Value = -69621

$ws.Range("H45").Value = 2035.1111
$ws.Range("I45").Value = 1584.5385
$ws.Range("K45").Value = 1584.5385
$ws.Range("M45").Value = -1207.5385

$ws.Range("H97").Value = 1286.7561
$ws.Range("I97").Value = 1035.5312
$ws.Range("K97").Value = 1035.5312
$ws.Range("M97").Value = -539.5311999999999

$ws.Range("H102").Value = 2988.2856
$ws.Range("I102").Value = 2874.25
$ws.Range("K102").Value = 2874.25
$ws.Range("M102").Value = -1252.25

$ws.Range("H110").Value = 645.2222
$ws.Range("I110").Value = 669.625
$ws.Range("J110").Value = 450
$ws.Range("K110").Value = 669.625
$ws.Range("L110").Value = 450
$ws.Range("M110").Value = 1375.375
$ws.Range("N110").Value = -4540

$ws.Range("H132").Value = 11239.186
$ws.Range("I132").Value = 6246.846
$ws.Range("K132").Value = 18740.538
$ws.Range("M132").Value = -16210.538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 299.8
$ws.Range("I22").Value = 299.8
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 299.8
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -126.8
$ws.Range("N22").ClearContents()

$ws.Range("H94").Value = 854.04346
$ws.Range("I94").Value = 358.375
$ws.Range("J94").Value = 1987
$ws.Range("K94").Value = 358.375
$ws.Range("L94").Value = 1987
$ws.Range("M94").Value = 92.625
$ws.Range("N94").Value = -2889

$ws.Range("H99").Value = 17060.143
$ws.Range("I99").Value = 34970
$ws.Range("K99").Value = 34970
$ws.Range("M99").Value = -33472

$ws.Range("H105").Value = 1267.2142
$ws.Range("I105").Value = 1133.0769
$ws.Range("K105").Value = 1133.0769
$ws.Range("M105").Value = 613.9231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2000
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H16").Value = 11933.286
$ws.Range("I16").Value = 13172.167
$ws.Range("K16").Value = 13172.167
$ws.Range("M16").Value = -12885.167

$ws.Range("H99").Value = 1979.1
$ws.Range("I99").Value = 1923.25
$ws.Range("K99").Value = 1923.25
$ws.Range("M99").Value = -425.25

$ws.Range("H113").Value = 11933.286
$ws.Range("I113").Value = 13172.167
$ws.Range("K113").Value = 13172.167
$ws.Range("M113").Value = -11002.167

$ws.Range("H126").Value = 1979.1
$ws.Range("I126").Value = 1923.25
$ws.Range("K126").Value = 5769.75
$ws.Range("M126").Value = -3299.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2349.3704
$ws.Range("J68").Value = 2238.7917
$ws.Range("L68").Value = 6716.375100000001
$ws.Range("N68").Value = -8338.375100000001

$ws.Range("H71").Value = 2349.3704
$ws.Range("J71").Value = 2238.7917
$ws.Range("L71").Value = 20149.1253
$ws.Range("N71").Value = -28261.1253

$ws.Range("H92").Value = 668053.25
$ws.Range("J92").Value = 830.6923
$ws.Range("L92").Value = 2492.0769
$ws.Range("N92").Value = -4988.0769

$ws.Range("H122").Value = 743.48
$ws.Range("I122").Value = 682
$ws.Range("K122").Value = 6138
$ws.Range("M122").Value = -3688

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2999.25
$ws.Range("I80").Value = 2999.5
$ws.Range("J80").Value = 2999
$ws.Range("K80").Value = 2999.5
$ws.Range("L80").Value = 2999
$ws.Range("M80").Value = -2001.5
$ws.Range("N80").Value = -4995

$ws.Range("H83").Value = 2999.25
$ws.Range("I83").Value = 2999.5
$ws.Range("J83").Value = 2999
$ws.Range("K83").Value = 14997.5
$ws.Range("L83").Value = 14995
$ws.Range("M83").Value = -10005.5
$ws.Range("N83").Value = -24979

$ws.Range("H102").Value = 3459.92
$ws.Range("I102").Value = 3035.8235
$ws.Range("J102").Value = 4361.125
$ws.Range("K102").Value = 3035.8235
$ws.Range("L102").Value = 4361.125
$ws.Range("M102").Value = -1413.8235
$ws.Range("N102").Value = -7605.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4236.727
$ws.Range("J40").Value = 4731.1665
$ws.Range("L40").Value = 4731.1665
$ws.Range("N40").Value = -5003.1665

$ws.Range("H46").Value = 2261.9375
$ws.Range("I46").Value = 2298.875
$ws.Range("J46").Value = 2225
$ws.Range("K46").Value = 2298.875
$ws.Range("L46").Value = 2225
$ws.Range("M46").Value = -2110.875
$ws.Range("N46").Value = -2601

$ws.Range("H122").Value = 6141.476
$ws.Range("I122").Value = 5468.8823
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 16406.6469
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -13956.6469
$ws.Range("N122").Value = -31900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 52765.8
$ws.Range("J105").Value = 52765.8
$ws.Range("L105").Value = 52765.8
$ws.Range("N105").Value = -59753.8

$ws.Range("H136").Value = 1425.9
$ws.Range("I136").Value = 1060.5385
$ws.Range("K136").Value = 3181.6155
$ws.Range("M136").Value = -631.6155000000003
